$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# New DRG-code rows for 'All Gastro' (T01A/T01B), Toxo/Typhoid (B06A, already
# present) and Listeria (T64B/T64C). Shared-string entries are created in
# the order the distinct text values are first written, so the write order
# below is chosen to reproduce that sequence.
# ----------------------------------------------------------------------

# 1) "T01A" first used here
$ws.Range("B42").Value = "T01A"

# 2) "DRG code" first used here (existing rows 38-41 previously had no
#    label in column D for the DRG-code rows; now they all get one)
$ws.Range("D38").Value = "DRG code"

# 3) "T01A/T01B"
$ws.Range("C44").Value = "T01A/T01B"

# 4) "Average of T01A and T01B"
$ws.Range("D44").Value = "Average of T01A and T01B"

# 5) "T01B"
$ws.Range("B43").Value = "T01B"

# 6) "T64B/T64C"
$ws.Range("C47").Value = "T64B/T64C"

# 7) "Average of T64B and T64C"
$ws.Range("D47").Value = "Average of T64B and T64C"

# 8) "T64B"
$ws.Range("B45").Value = "T64B"

# 9) "T64C"
$ws.Range("B46").Value = "T64C"

# ----------------------------------------------------------------------
# Remaining writes that re-use the shared strings created above.
# ----------------------------------------------------------------------
$ws.Range("C42").Value = "T01A"
$ws.Range("C43").Value = "T01B"
$ws.Range("C45").Value = "T64B"
$ws.Range("C46").Value = "T64C"

$ws.Range("D39").Value = "DRG code"
$ws.Range("D40").Value = "DRG code"
$ws.Range("D41").Value = "DRG code"
$ws.Range("D42").Value = "DRG code"
$ws.Range("D43").Value = "DRG code"

# ----------------------------------------------------------------------
# Numeric values / formulas (do not touch the shared-string table).
# ----------------------------------------------------------------------
$ws.Range("A42").Value = 50828
$ws.Range("A43").Value = 18592
$ws.Range("A44").Formula = "=AVERAGE(A42:A43)"

$ws.Range("A45").Value = 12233
$ws.Range("A46").Value = 5199
$ws.Range("A47").Formula = "=AVERAGE(A45:A46)"

# ----------------------------------------------------------------------
# Formatting: copy number format / font from the matching, already-styled
# cells directly above the new rows, so styles line up with the rest of
# the column.
# ----------------------------------------------------------------------
$ws.Range("A37").Copy()
$ws.Range("A42:A47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B37").Copy()
$ws.Range("B42:C43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B45:C46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("C37").Copy()
$ws.Range("C44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("D37").Copy()
$ws.Range("D39:D44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# View / selection updates (window scrolled down a few rows, selection
# moved to the newly added last cell).
# ----------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("A47").Select()
